$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "naturedossiers_id" (column L) values for the affected rows
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("L11").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("L13").Value = 1
$ws.Range("L14").Value = 1
$ws.Range("L15").Value = 1
$ws.Range("L25").Value = 2
$ws.Range("L26").Value = 2
$ws.Range("L27").Value = 2
$ws.Range("L28").Value = 2

# Move the active selection, matching where the author last clicked
$ws.Range("M31").Select()
